$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2000
$ws.Range("J48").Value = 2000
$ws.Range("L48").Value = 6000
$ws.Range("N48").Value = -6584
$ws.Range("H56").Value = 2000
$ws.Range("J56").Value = 2000
$ws.Range("L56").Value = 6000
$ws.Range("N56").Value = -7068
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 3486.5264
$ws.Range("I113").Value = 2260
$ws.Range("J113").Value = 3924.5715
$ws.Range("K113").Value = 2260
$ws.Range("L113").Value = 3924.5715
$ws.Range("M113").Value = 994
$ws.Range("N113").Value = -10432.5715
$ws.Range("H138").Value = 1885.3256
$ws.Range("I138").Value = 1255.0667
$ws.Range("J138").Value = 3339.7693
$ws.Range("K138").Value = 3765.2001
$ws.Range("L138").Value = 10019.3079
$ws.Range("M138").Value = 1374.7999
$ws.Range("N138").Value = -20299.3079

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3210.2354
$ws.Range("I2").Value = 2133.111
$ws.Range("J2").Value = 4422
$ws.Range("K2").Value = 2133.111
$ws.Range("L2").Value = 4422
$ws.Range("M2").Value = -2020.111
$ws.Range("N2").Value = -4648
$ws.Range("H32").Value = 5437.075
$ws.Range("I32").Value = 5845.154
$ws.Range("J32").Value = 3315.0667
$ws.Range("K32").Value = 5845.154
$ws.Range("L32").Value = 3315.0667
$ws.Range("M32").Value = -5558.154
$ws.Range("N32").Value = -3889.0667
$ws.Range("H45").Value = 6437.4
$ws.Range("I45").Value = 3800
$ws.Range("J45").Value = 7096.75
$ws.Range("K45").Value = 3800
$ws.Range("L45").Value = 7096.75
$ws.Range("M45").Value = -3423
$ws.Range("N45").Value = -7850.75
$ws.Range("H60").Value = 3364
$ws.Range("I60").Value = 3364
$ws.Range("K60").Value = 3364
$ws.Range("M60").Value = -2631
$ws.Range("H88").Value = 1887.4445
$ws.Range("J88").Value = 1855.2858
$ws.Range("L88").Value = 1855.2858
$ws.Range("N88").Value = -2667.2858
$ws.Range("H91").Value = 1887.4445
$ws.Range("J91").Value = 1855.2858
$ws.Range("L91").Value = 1855.2858
$ws.Range("N91").Value = -4663.2858
$ws.Range("H110").Value = 2424.158
$ws.Range("I110").Value = 1057.2667
$ws.Range("J110").Value = 7550
$ws.Range("K110").Value = 1057.2667
$ws.Range("L110").Value = 7550
$ws.Range("M110").Value = 987.7333000000001
$ws.Range("N110").Value = -11640
$ws.Range("H116").Value = 3210.2354
$ws.Range("I116").Value = 2133.111
$ws.Range("J116").Value = 4422
$ws.Range("K116").Value = 2133.111
$ws.Range("L116").Value = 4422
$ws.Range("M116").Value = 160.8890000000001
$ws.Range("N116").Value = -9010
$ws.Range("H132").Value = 4862.3423
$ws.Range("I132").Value = 1791.6471
$ws.Range("J132").Value = 7348.143
$ws.Range("K132").Value = 5374.9413
$ws.Range("L132").Value = 22044.429
$ws.Range("M132").Value = -2844.9413
$ws.Range("N132").Value = -27104.429

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3210.2354
$ws.Range("I3").Value = 2133.111
$ws.Range("J3").Value = 4422
$ws.Range("K3").Value = 2133.111
$ws.Range("L3").Value = 4422
$ws.Range("M3").Value = -2019.111
$ws.Range("N3").Value = -4650
$ws.Range("H86").Value = 1264
$ws.Range("I86").Value = 1397.8889
$ws.Range("J86").Value = 1177.9286
$ws.Range("K86").Value = 1397.8889
$ws.Range("L86").Value = 1177.9286
$ws.Range("M86").Value = -274.8888999999999
$ws.Range("N86").Value = -3423.9286
$ws.Range("H89").Value = 1264
$ws.Range("I89").Value = 1397.8889
$ws.Range("J89").Value = 1177.9286
$ws.Range("K89").Value = 6989.4445
$ws.Range("L89").Value = 5889.643
$ws.Range("M89").Value = -1373.4445
$ws.Range("N89").Value = -17121.643
$ws.Range("H134").Value = 8137.963
$ws.Range("I134").Value = 5419.5
$ws.Range("J134").Value = 8914.666999999999
$ws.Range("K134").Value = 16258.5
$ws.Range("L134").Value = 26744.001
$ws.Range("M134").Value = -13723.5
$ws.Range("N134").Value = -31814.001

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4452.4165
$ws.Range("I7").Value = 8356.583000000001
$ws.Range("J7").Value = 548.25
$ws.Range("K7").Value = 8356.583000000001
$ws.Range("L7").Value = 548.25
$ws.Range("M7").Value = -8243.583000000001
$ws.Range("N7").Value = -774.25
$ws.Range("H31").Value = 6174858.5
$ws.Range("I31").Value = 1608.8718
$ws.Range("J31").Value = 22225308
$ws.Range("K31").Value = 1608.8718
$ws.Range("L31").Value = 22225308
$ws.Range("M31").Value = -1313.8718
$ws.Range("N31").Value = -22225898
$ws.Range("H34").Value = 6174858.5
$ws.Range("I34").Value = 1608.8718
$ws.Range("J34").Value = 22225308
$ws.Range("K34").Value = 1608.8718
$ws.Range("L34").Value = 22225308
$ws.Range("M34").Value = -1406.8718
$ws.Range("N34").Value = -22225712
$ws.Range("H58").Value = 1354305.1
$ws.Range("I58").Value = 2001.6957
$ws.Range("J58").Value = 3575946.2
$ws.Range("K58").Value = 2001.6957
$ws.Range("L58").Value = 3575946.2
$ws.Range("M58").Value = -1798.6957
$ws.Range("N58").Value = -3576352.2
$ws.Range("H63").Value = 23168.1
$ws.Range("J63").Value = 23168.1
$ws.Range("L63").Value = 23168.1
$ws.Range("N63").Value = -24540.1
$ws.Range("H66").Value = 23168.1
$ws.Range("J66").Value = 23168.1
$ws.Range("L66").Value = 69504.29999999999
$ws.Range("N66").Value = -76368.29999999999
$ws.Range("H100").Value = 39890
$ws.Range("J100").Value = 39890
$ws.Range("L100").Value = 39890
$ws.Range("N100").Value = -42054
$ws.Range("H107").Value = 395.5
$ws.Range("I107").Value = 306.72726
$ws.Range("J107").Value = 721
$ws.Range("K107").Value = 306.72726
$ws.Range("L107").Value = 721
$ws.Range("M107").Value = 1613.27274
$ws.Range("N107").Value = -4561
$ws.Range("H132").Value = 2481.625
$ws.Range("I132").Value = 1579.0667
$ws.Range("J132").Value = 3278
$ws.Range("K132").Value = 4737.2001
$ws.Range("L132").Value = 9834
$ws.Range("M132").Value = -2207.2001
$ws.Range("N132").Value = -14894
$ws.Range("H134").Value = 2831.5652
$ws.Range("I134").Value = 1584.3636
$ws.Range("J134").Value = 3974.8333
$ws.Range("K134").Value = 4753.0908
$ws.Range("L134").Value = 11924.4999
$ws.Range("M134").Value = -2218.0908
$ws.Range("N134").Value = -16994.4999
$ws.Range("H136").Value = 1354305.1
$ws.Range("I136").Value = 2001.6957
$ws.Range("J136").Value = 3575946.2
$ws.Range("K136").Value = 6005.0871
$ws.Range("L136").Value = 10727838.6
$ws.Range("M136").Value = -3455.0871
$ws.Range("N136").Value = -10732938.6

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1121.5
$ws.Range("I5").Value = 584.5714
$ws.Range("K5").Value = 1753.7142
$ws.Range("M5").Value = -1641.7142
$ws.Range("H7").Value = 600
$ws.Range("J7").Value = 766.6667
$ws.Range("L7").Value = 2300.0001
$ws.Range("N7").Value = -2524.0001
$ws.Range("H68").Value = 1665.6666
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1665.6666
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4996.9998
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6618.9998
$ws.Range("H71").Value = 1665.6666
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1665.6666
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14990.9994
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -23102.9994
$ws.Range("H95").Value = 6333.3335
$ws.Range("J95").Value = 6333.3335
$ws.Range("L95").Value = 19000.0005
$ws.Range("N95").Value = -23118.0005
$ws.Range("H114").Value = 1145.5454
$ws.Range("I114").Value = 402.25
$ws.Range("J114").Value = 1570.2858
$ws.Range("K114").Value = 1206.75
$ws.Range("L114").Value = 4710.857400000001
$ws.Range("M114").Value = 2047.25
$ws.Range("N114").Value = -11218.8574
$ws.Range("H122").Value = 2710.8293
$ws.Range("I122").Value = 448.875
$ws.Range("J122").Value = 3259.182
$ws.Range("K122").Value = 4039.875
$ws.Range("L122").Value = 29332.638
$ws.Range("M122").Value = -1589.875
$ws.Range("N122").Value = -34232.638
$ws.Range("H129").Value = 2205.72
$ws.Range("I129").Value = 948.3333
$ws.Range("J129").Value = 2602.7896
$ws.Range("K129").Value = 2844.9999
$ws.Range("L129").Value = 7808.3688
$ws.Range("M129").Value = 2155.0001
$ws.Range("N129").Value = -17808.3688
$ws.Range("H131").Value = 205234.4
$ws.Range("I131").Value = 2000396
$ws.Range("J131").Value = 1238.7727
$ws.Range("K131").Value = 6001188
$ws.Range("L131").Value = 3716.3181
$ws.Range("M131").Value = -5996148
$ws.Range("N131").Value = -13796.3181
$ws.Range("H135").Value = 1121.5
$ws.Range("I135").Value = 584.5714
$ws.Range("K135").Value = 5261.1426
$ws.Range("M135").Value = -2726.1426

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H86").Value = 20000
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18814
$ws.Range("H89").Value = 20000
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 60000
$ws.Range("M89").Value = -54072
$ws.Range("H132").Value = 52635772
$ws.Range("I132").Value = 111115960
$ws.Range("J132").Value = 3599.8
$ws.Range("K132").Value = 333347880
$ws.Range("L132").Value = 10799.4
$ws.Range("M132").Value = -333345350
$ws.Range("N132").Value = -15859.4
$ws.Range("H136").Value = 33335746
$ws.Range("I136").Value = 45455544
$ws.Range("J136").Value = 6300
$ws.Range("K136").Value = 136366632
$ws.Range("L136").Value = 18900
$ws.Range("M136").Value = -136364082
$ws.Range("N136").Value = -24000

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 880.8077
$ws.Range("I113").Value = 481.2857
$ws.Range("J113").Value = 1028
$ws.Range("K113").Value = 1443.8571
$ws.Range("L113").Value = 3084
$ws.Range("M113").Value = 726.1428999999998
$ws.Range("N113").Value = -7424
$ws.Range("H125").Value = 49800
$ws.Range("J125").Value = 49800
$ws.Range("L125").Value = 49800
$ws.Range("N125").Value = -59640
$ws.Range("H133").Value = 42715
$ws.Range("J133").Value = 42715
$ws.Range("L133").Value = 42715
$ws.Range("N133").Value = -52835
$ws.Range("H136").Value = 7144720
$ws.Range("I136").Value = 16668120
$ws.Range("J136").Value = 2170.15
$ws.Range("K136").Value = 50004360
$ws.Range("L136").Value = 6510.450000000001
$ws.Range("M136").Value = -50001810
$ws.Range("N136").Value = -11610.45
